$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 382, shifting existing rows 382-471 down to 383-472.
$ws.Rows("382:382").Insert()

# Populate the newly inserted row 382 with the new record.
$ws.Range("A382").Value = 3
$ws.Range("B382").Value = "Femacal de La Calera"
$ws.Range("C382").Value = "Coquimbo"
$ws.Range("D382").Value = 44889
$ws.Range("E382").Value = 5
$ws.Range("F382").Value = 100112031
$ws.Range("G382").Value = "Poroto verde"
$ws.Range("H382").Value = "Magnum"
$ws.Range("I382").Value = "Primera"
$ws.Range("J382").Value = 85
$ws.Range("K382").Value = 40000
$ws.Range("L382").Value = 41000
$ws.Range("M382").Value = 40471
$ws.Range("N382").Value = "`$/malla 25 kilos"
$ws.Range("O382").Value = "Provincia de Quillota"
$ws.Range("P382").Value = 1619
$ws.Range("Q382").Value = 25
$ws.Range("R382").Value = "Hortaliza"
